$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRows = @(
    @("WV50 FILTER", "Fallo atornillador", "2024-06-10", "10:56:38", "Mañana", "10:56:40", "0:00:02", "-0.00 minutos"),
    @("WV50 FILTER", "Fallo cámara cover", "2024-06-10", "10:56:51", "Mañana", "10:56:56", "0:00:05", "0.05 minutos"),
    @("WV50 FILTER", "Traza", "2024-06-10", "10:57:07", "Mañana", "10:57:11", "0:00:04", "0.08 minutos"),
    @("WV50 FILTER", "Robot no coloca bien filter en palet", "2024-06-10", "10:58:13", "Mañana", "10:58:18", "0:00:05", "0.20 minutos"),
    @("WV50 FILTER", "NOK Soldadura Plástico", "2024-06-10", "10:58:32", "Mañana", "10:58:33", "0:00:01", "0.19 minutos")
)

$startRow = 166
$endRow = $startRow + $newRows.Count - 1

# Force text format on the "Fecha" column so Excel does not auto-convert
# the date-looking string (e.g. "2024-06-10") into a date serial number.
$ws.Range("C$startRow`:C$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
